$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KTOS")
$ws.Columns("D").Insert()
